$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.635.55'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.473.85'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.14'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.36'
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.551'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0867'
$ws.Range("E10").Value = '  +9.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '33.15'
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '2.856.47'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.54'
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("D16").Value = '2.489.61'
$ws.Range("E16").Value = '  +1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.791'
$ws.Range("E17").Value = '  +2.98%  '
$ws.Range("D18").Value = '41.598.93'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.45'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.74'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.28'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.09'
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  +1.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.85'
$ws.Range("E27").Value = '  +2.65%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.69'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.57'
$ws.Range("E30").Value = '  +4.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.88'
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0767'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.30'
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.02'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("E42").Value = '  +3.25%  '
$ws.Range("D43").Value = '1.987.00'
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.75'
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.97'
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.46'
$ws.Range("E47").Value = '  +5.94%  '
$ws.Range("D48").Value = '2.714.33'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.95'
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.74'
$ws.Range("E50").Value = '  +5.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.08'

